$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "3.2. Thay đổi trạng thái nhiều sản phẩm"
$ws.Range("C4").Value = "https://github.com/nguyentienminh07102004/product-management/commit/e896019a677ede65944dae61c70a40c4ecc8a67d"

$ws.Range("C4").Select()
